$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Copy of lexer-FSM"

# Align S_BLANK (row 3) transitions to S_START (row 2) transitions
$ws.Range("D2:AF2").Copy($ws.Range("D3:AF3"))

# Closing-block transition from S_BLANK should be T_BLK_OP (like the opening
# block column), not T_BLK_CL as copied from S_START
$ws.Range("G3").Value2 = $ws.Range("F3").Value2

# Update the remembered selection on the sheet
$null = $ws.Range("AF3").Select()
